$d = $word.ActiveDocument

$replacements = @(
    @{old = "990÷8=123, 6"; new = "472÷9=52, 4"},
    @{old = "214÷3=71, 1"; new = "943÷7=134, 5"},
    @{old = "368÷9=40, 8"; new = "250÷7=35, 5"},
    @{old = "512÷2=256, 0"; new = "555÷6=92, 3"},
    @{old = "708÷2=354, 0"; new = "431÷5=86, 1"},
    @{old = "427÷4=106, 3"; new = "282÷8=35, 2"},
    @{old = "397÷4=99, 1"; new = "600÷7=85, 5"},
    @{old = "840÷3=280, 0"; new = "417÷8=52, 1"},
    @{old = "511÷7=73, 0"; new = "787÷5=157, 2"},
    @{old = "240÷6=40, 0"; new = "637÷2=318, 1"},
    @{old = "372÷7=53, 1"; new = "949÷4=237, 1"},
    @{old = "548÷8=68, 4"; new = "697÷9=77, 4"},
    @{old = "538÷7=76, 6"; new = "642÷6=107, 0"},
    @{old = "810÷8=101, 2"; new = "668÷8=83, 4"},
    @{old = "224÷5=44, 4"; new = "519÷9=57, 6"},
    @{old = "366÷9=40, 6"; new = "451÷8=56, 3"},
    @{old = "236÷9=26, 2"; new = "983÷6=163, 5"},
    @{old = "847÷2=423, 1"; new = "577÷8=72, 1"},
    @{old = "194÷2=97, 0"; new = "319÷6=53, 1"},
    @{old = "439÷4=109, 3"; new = "388÷9=43, 1"},
    @{old = "776÷7=110, 6"; new = "916÷7=130, 6"},
    @{old = "562÷4=140, 2"; new = "196÷3=65, 1"},
    @{old = "318÷7=45, 3"; new = "755÷4=188, 3"},
    @{old = "186÷4=46, 2"; new = "493÷2=246, 1"},
    @{old = "763÷5=152, 3"; new = "405÷4=101, 1"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}

$d.Save()
